$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H341").Value = 'Кроличья лапка'
$ws.Range("H342").Value = 'Некоторые суеверные люди думают что это амулет, способный принести удачу.'
$ws.Range("H344").Value = 'Крепящееся к пальцам дробящее оружие с выпуклостью или шипами направленными наружу. Намного мощнее обычного кулака.'
$ws.Range("H345").Value = 'Обычно используется для игры в бейсбол, но его сподручность позволяет бить не только мячи.'
$ws.Range("H346").Value = 'Бейсбольная бита'
$ws.Range("H347").Value = 'Кастет'
$ws.Range("H348").Value = 'Мачете'
$ws.Range("H349").Value = 'Длинное лезвие обычно используеещеся для колки древесины, а иногда и для копания.'
$ws.Range("H350").Value = 'Старая рукопись'
$ws.Range("H352").Value = 'ВОЛ+30'
$ws.Range("H353").Value = 'Медицина+20'
$ws.Range("H354").Value = 'ТЕЛ+30'
$ws.Range("H355").Value = 'СИЛ+30'
$ws.Range("H356").Value = 'Скрытность+10'
$ws.Range("H357").Value = 'Экипировка: HP+3 каждый ход в бою'
$ws.Range("H358").Value = 'Броня+3'
$ws.Range("H359").Value = 'Спорт+50'
$ws.Range("H360").Value = 'Наблюдательность+20'
$ws.Range("H361").Value = 'Дипломатия+10'
$ws.Range("H362").Value = 'Эрудиция+20'
$ws.Range("H364").Value = 'Чтение'
$ws.Range("H366").Value = 'Вырезка из газеты, на ней объявление о пропаже пианиста Паскаля из Таунсенда.'
$ws.Range("H367").Value = 'Вырезка из новостей'
$ws.Range("H368").Value = 'Ключ от кладовой'
$ws.Range("H369").Value = 'Ключ от камеры хранения на втором этаже отеля.'
$ws.Range("H370").Value = 'Ключ Райана'
$ws.Range("H380").Value = 'Электромагнитная граната'
$ws.Range("H398").Value = 'Виноградная гирлянда'
$ws.Range("H399").Value = 'Браслет из виноградных лоз с нераспустившимся цветком на нем.'
$ws.Range("H400").Value = 'Восстановить руну'
$ws.Range("H401").Value = 'Спирт'
$ws.Range("H402").Value = 'Руна защиты'
$ws.Range("H403").Value = 'Руна души'
$ws.Range("H404").Value = 'Руна разума'
$ws.Range("H405").Value = 'Руна скорости'
$ws.Range("H406").Value = 'Снаряд производит высокочастотные электромагнитные волны, нарушающие нормальные физиологические функции организмов.'
$ws.Range("H407").Value = 'Обычное оружие уличных протестующих.'
$ws.Range("H452").Value = 'Если на поле только один противник, атака получит бонусный кубик.'
$ws.Range("H510").Value = 'Экипировка:
Получение урона уменьшит прочность на 1.'
$ws.Range("H511").Value = 'Экипировка: ЛОВ-20
Получение урона уменьшит прочность на 1.'
$ws.Range("H512").Value = 'Экипировка: Скорость+50
Каждый ход в бою прочность этого предмета уменьшится на 1.'
$ws.Range("H513").Value = 'Полностью восстанавливает MP.'
$ws.Range("H514").Value = 'Каждый полученный урон уменьшает прочность на 1. Будет уничтожен когда прочность достигнет 0.'
$ws.Range("H515").Value = 'Каждый полученный урон уменьшает прочность на 1.'
$ws.Range("H516").Value = 'В начале хода восставливает 3 HP и теряет 1 прочность.'
$ws.Range("H517").Value = 'Пачка боеприпасов'
$ws.Range("H518").Value = 'SP+[1-3]
Может использоваться только если SP < 50'
$ws.Range("H519").Value = 'SP+[1-4], снимает [Замешательство]
Может использоваться только если SP < 65'
$ws.Range("H520").Value = 'Экипировка: После каждой атаки теряет 1 единицу прочности. Не ломается.'
$ws.Range("H521").Value = 'Экипировка: При атаке потеряет 1 прочность вместо того чтобы расходовать боеприпасы.
Не ломается.'
$ws.Range("H522").Value = 'Деревянный лук'
$ws.Range("H523").Value = 'Мощное оружие дальнего боя. Для постоянного использования требуется большая физическая сила.'
$ws.Range("H524").Value = 'Палочка-дразнилка для кошек'
$ws.Range("H525").Value = 'Это мост дружбы соеденяющий вас с плюшевыми милашками!'
$ws.Range("H526").Value = 'Экипировка: Дипломатия +5, SP +5
Если ваша цель - котенок, ваша проверка дипломатии увенчается критическим успехом!'
$ws.Range("H527").Value = 'Этот предмет нельзя снимать до окончания модуля. Вы хотите его надеть?'
$ws.Range("H528").Value = 'Ключ управления колесом обозрения'
$ws.Range("H529").Value = 'Это ключ к переключателю управления колесом обозрения'
$ws.Range("H530").Value = 'Липкая штука'
$ws.Range("H531").Value = 'Случайный продукт вашей обработки товара, хотя сама обработка не удалась, но вещь выглядит так, как будто она съедобна'